# Update Mappings 22 Ontologies
# Adds a new "SBO_DEF" column (F) with header + value, and updates the
# existing SBO_IRI annotation mapping description (C2) to the richer
# dictionary representation produced by the newer ontology mapper.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell F1: "SBO_DEF" ---------------------------------------
# Copy formatting from the existing header cell E1 (bold font + border)
# before writing the new header text, so F1 matches the other headers.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F1").Value = "SBO_DEF"

# --- Updated mapping text for C2 ------------------------------------------
$ws.Range("C2").Value = "{'label': None, 'prefLabel': None, 'altLabel': None, 'name': 'annotation'}"

# --- New data cell F2: "[]" ------------------------------------------------
$ws.Range("F2").Value = "[]"

$ws.Application.CutCopyMode = $false
